# Weekly update to the "Femacal de La Calera - Achicoria" price sheet.
# A new daily record is inserted at row 57 (pushing the existing rows
# 57:159 down to 58:160); the rest of the data is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 57 - this shifts rows
# 57:159 down to 58:160 and grows the sheet's used range to A1:R160.
$ws.Rows.Item(57).Insert()

# Fill in the new record's values.
$ws.Cells.Item(57, 1).Value = 3
$ws.Cells.Item(57, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(57, 3).Value = "Coquimbo"
$ws.Cells.Item(57, 4).Value = 44533
$ws.Cells.Item(57, 5).Value = 5
$ws.Cells.Item(57, 6).Value = 100112010
$ws.Cells.Item(57, 7).Value = "Achicoria"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 60
$ws.Cells.Item(57, 11).Value = 5500
$ws.Cells.Item(57, 12).Value = 5500
$ws.Cells.Item(57, 13).Value = 5500
$ws.Cells.Item(57, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(57, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(57, 16).Value = 344
$ws.Cells.Item(57, 17).Value = 16
$ws.Cells.Item(57, 18).Value = "Hortaliza"
